$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.1719888661583355
$ws.Range("D2").Value = 0.8650185693849215

$ws.Range("C3").Value = 0.02483749757979744
$ws.Range("D3").Value = 0.9804084972311169

$ws.Range("C4").Value = 0.3072838632568332
$ws.Range("D4").Value = 0.7615172468692648

$ws.Range("C5").Value = -0.410206562334672
$ws.Range("D5").Value = 0.6856217672391547

$ws.Range("C6").Value = 0.2130661159223731
$ws.Range("D6").Value = 0.8332374139446477

$ws.Range("C7").Value = 0.7022547868001311
$ws.Range("D7").Value = 0.4898843377065685

$ws.Range("C8").Value = -0.3310935595467661
$ws.Range("D8").Value = 0.7437067712476315

$ws.Range("C9").Value = 0.2771722403643972
$ws.Range("D9").Value = 0.7842355650058841

$ws.Range("C10").Value = -0.4875399770126377
$ws.Range("D10").Value = 0.630698778436797

$ws.Range("C11").Value = -0.7079446970155765
$ws.Range("D11").Value = 0.4864129525854146
